$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (E1:G1) ---
$ws.Range("E1").Value = "descricao_longa"
$ws.Range("F1").Value = "especificacoes"
$ws.Range("G1").Value = "estoque"

# --- Column B: replace product codes with descriptive names ---
$ws.Range("B2").Value = "Corpo de Borboleta Fiat Palio"
$ws.Range("B3").Value = "Corpo de Borboleta Fiat Uno/Mobi"
$ws.Range("B4").Value = "Corpo de Borboleta Fiat Linea"
$ws.Range("B5").Value = "Corpo de Borboleta GM Celta"
$ws.Range("B6").Value = "Corpo de Borboleta GM Cruze 1.4"

# --- Column C: update row 3 detail text ---
$ws.Range("C3").Value = "Uno Celebration / Mobi / Palio 1.0"

# --- Column E: descricao_longa ---
$ws.Range("E2").Value = "Corpo de borboleta completo, com tampa e asas, abre e fecha e foi feita da mais refinada IA. Serve no Palio e vai na Strada tbm"
$ws.Range("E3").Value = "Corpinho para o Uno e Mobi, é pequeno igual eles "
$ws.Range("E4").Value = "Oxi, existe carro 1.9? Melhor fazer 2.0 logo de uma vez po"
$ws.Range("E5").Value = "Mais um corpo de borboleta "
$ws.Range("E6").Value = "Encorpado na borboleta "

# --- Column F: especificacoes (only rows 2 and 3) ---
$ws.Range("F2").Value = "Motores 1.4"
$ws.Range("F3").Value = "Motores 1.0"

# --- Column G: estoque ---
$ws.Range("G2").Value = 10
$ws.Range("G3").Value = 12
$ws.Range("G4").Value = 7
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 2

# --- Match body formatting (vertical-center) used by columns A-D on the new E/F data cells ---
$ws.Range("A2").Copy()
$ws.Range("E2:F3").PasteSpecial(-4122)
$ws.Range("E4:E6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths for the new columns (target widths: 15.28515625 / 14 / 8.28515625) ---
$ws.Columns.Item(5).ColumnWidth = 14.451822916666666
$ws.Columns.Item(6).ColumnWidth = 13.166666666666666
$ws.Columns.Item(7).ColumnWidth = 7.451822916666667

# --- Selection matches the final saved state ---
$ws.Range("G9").Select()
